$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5358
$ws.Range("I32").Value = 3950
$ws.Range("K32").Value = 3950
$ws.Range("M32").Value = -3624
$ws.Range("H43").Value = 5419
$ws.Range("J43").Value = 7484.25
$ws.Range("L43").Value = 7484.25
$ws.Range("N43").Value = -7622.25
$ws.Range("H76").Value = 4353.5264
$ws.Range("I76").Value = 3825.6667
$ws.Range("J76").Value = 5258.4287
$ws.Range("K76").Value = 3825.6667
$ws.Range("L76").Value = 5258.4287
$ws.Range("M76").Value = -3510.6667
$ws.Range("N76").Value = -5888.4287
$ws.Range("H79").Value = 4353.5264
$ws.Range("I79").Value = 3825.6667
$ws.Range("J79").Value = 5258.4287
$ws.Range("K79").Value = 3825.6667
$ws.Range("L79").Value = 5258.4287
$ws.Range("M79").Value = -2733.6667
$ws.Range("N79").Value = -7442.4287
$ws.Range("H137").Value = 6595.0244
$ws.Range("I137").Value = 8629.448
$ws.Range("K137").Value = 25888.344
$ws.Range("M137").Value = -23338.344
$ws.Range("H138").Value = 3799.4363
$ws.Range("J138").Value = 4428.6523
$ws.Range("L138").Value = 13285.9569
$ws.Range("N138").Value = -23565.9569

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9359.706
$ws.Range("I2").Value = 10276.071
$ws.Range("K2").Value = 10276.071
$ws.Range("M2").Value = -10163.071
$ws.Range("H5").Value = 2434.1428
$ws.Range("I5").Value = 2731.5
$ws.Range("K5").Value = 2731.5
$ws.Range("M5").Value = -2619.5
$ws.Range("H37").Value = 5966
$ws.Range("I37").Value = 5966
$ws.Range("K37").Value = 5966
$ws.Range("M37").Value = -5693
$ws.Range("H63").Value = 5835.6665
$ws.Range("I63").Value = 6667.6665
$ws.Range("K63").Value = 6667.6665
$ws.Range("M63").Value = -5981.6665
$ws.Range("H66").Value = 5835.6665
$ws.Range("I66").Value = 6667.6665
$ws.Range("K66").Value = 33338.3325
$ws.Range("M66").Value = -29906.3325
$ws.Range("H74").Value = 4105.41
$ws.Range("I74").Value = 2620.3667
$ws.Range("K74").Value = 2620.3667
$ws.Range("M74").Value = -1746.3667
$ws.Range("H77").Value = 4105.41
$ws.Range("I77").Value = 2620.3667
$ws.Range("K77").Value = 13101.8335
$ws.Range("M77").Value = -8733.833500000001
$ws.Range("H88").Value = 3199.7144
$ws.Range("I88").Value = 2799.6667
$ws.Range("J88").Value = 3499.75
$ws.Range("K88").Value = 2799.6667
$ws.Range("L88").Value = 3499.75
$ws.Range("M88").Value = -2393.6667
$ws.Range("N88").Value = -4311.75
$ws.Range("H91").Value = 3199.7144
$ws.Range("I91").Value = 2799.6667
$ws.Range("J91").Value = 3499.75
$ws.Range("K91").Value = 2799.6667
$ws.Range("L91").Value = 3499.75
$ws.Range("M91").Value = -1395.6667
$ws.Range("N91").Value = -6307.75
$ws.Range("H116").Value = 9359.706
$ws.Range("I116").Value = 10276.071
$ws.Range("K116").Value = 10276.071
$ws.Range("M116").Value = -7982.071

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9359.706
$ws.Range("I3").Value = 10276.071
$ws.Range("K3").Value = 10276.071
$ws.Range("M3").Value = -10162.071
$ws.Range("H4").Value = 2434.1428
$ws.Range("I4").Value = 2731.5
$ws.Range("K4").Value = 2731.5
$ws.Range("M4").Value = -2616.5
$ws.Range("H20").Value = 6080.3335
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 6080.3335
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 6080.3335
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -6574.3335
$ws.Range("H86").Value = 24228.572
$ws.Range("I86").Value = 15920
$ws.Range("K86").Value = 15920
$ws.Range("M86").Value = -14797
$ws.Range("H89").Value = 24228.572
$ws.Range("I89").Value = 15920
$ws.Range("K89").Value = 79600
$ws.Range("M89").Value = -73984
$ws.Range("H105").Value = 4877.3335
$ws.Range("J105").Value = 6166.3335
$ws.Range("L105").Value = 6166.3335
$ws.Range("N105").Value = -9660.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3257.303
$ws.Range("I31").Value = 2618.647
$ws.Range("K31").Value = 2618.647
$ws.Range("M31").Value = -2323.647
$ws.Range("H34").Value = 3257.303
$ws.Range("I34").Value = 2618.647
$ws.Range("K34").Value = 2618.647
$ws.Range("M34").Value = -2416.647
$ws.Range("H58").Value = 1840.5217
$ws.Range("I58").Value = 1572.7142
$ws.Range("J58").Value = 2257.111
$ws.Range("K58").Value = 1572.7142
$ws.Range("L58").Value = 2257.111
$ws.Range("M58").Value = -1369.7142
$ws.Range("N58").Value = -2663.111
$ws.Range("H107").Value = 22973.6
$ws.Range("I107").Value = 37035
$ws.Range("J107").Value = 1881.5
$ws.Range("K107").Value = 37035
$ws.Range("L107").Value = 1881.5
$ws.Range("M107").Value = -35115
$ws.Range("N107").Value = -5721.5
$ws.Range("H136").Value = 1840.5217
$ws.Range("I136").Value = 1572.7142
$ws.Range("J136").Value = 2257.111
$ws.Range("K136").Value = 4718.142599999999
$ws.Range("L136").Value = 6771.333
$ws.Range("M136").Value = -2168.142599999999
$ws.Range("N136").Value = -11871.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 897.1
$ws.Range("I50").Value = 492.42856
$ws.Range("K50").Value = 1477.28568
$ws.Range("M50").Value = -996.28568
$ws.Range("H51").Value = 2561.75
$ws.Range("J51").Value = 2785.1428
$ws.Range("L51").Value = 8355.428400000001
$ws.Range("N51").Value = -9275.428400000001
$ws.Range("H53").Value = 897.1
$ws.Range("I53").Value = 492.42856
$ws.Range("K53").Value = 1477.28568
$ws.Range("M53").Value = -996.28568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 510833.16
$ws.Range("I29").Value = 1006666.3
$ws.Range("K29").Value = 1006666.3
$ws.Range("M29").Value = -1006376.3
$ws.Range("H38").Value = 17500
$ws.Range("J38").Value = 17500
$ws.Range("L38").Value = 17500
$ws.Range("N38").Value = -18426
$ws.Range("H70").Value = 6703.737
$ws.Range("I70").Value = 6118.778
$ws.Range("K70").Value = 6118.778
$ws.Range("M70").Value = -5848.778
$ws.Range("H73").Value = 6703.737
$ws.Range("I73").Value = 6118.778
$ws.Range("K73").Value = 6118.778
$ws.Range("M73").Value = -5182.778
$ws.Range("H107").Value = 546.4211
$ws.Range("I107").Value = 817
$ws.Range("J107").Value = 302.9
$ws.Range("K107").Value = 817
$ws.Range("L107").Value = 302.9
$ws.Range("M107").Value = 1103
$ws.Range("N107").Value = -4142.9
$ws.Range("H122").Value = 33089.8
$ws.Range("I122").Value = 26733
$ws.Range("K122").Value = 80199
$ws.Range("M122").Value = -77749
$ws.Range("H126").Value = 22604.422
$ws.Range("I126").Value = 42874.75
$ws.Range("K126").Value = 128624.25
$ws.Range("M126").Value = -126154.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 430.22223
$ws.Range("I9").Value = 403.14285
$ws.Range("J9").Value = 525
$ws.Range("K9").Value = 403.14285
$ws.Range("L9").Value = 525
$ws.Range("M9").Value = -179.14285
$ws.Range("N9").Value = -973
$ws.Range("H46").Value = 2286.611
$ws.Range("I46").Value = 2509.1428
$ws.Range("J46").Value = 2145
$ws.Range("K46").Value = 2509.1428
$ws.Range("L46").Value = 2145
$ws.Range("M46").Value = -2321.1428
$ws.Range("N46").Value = -2521
$ws.Range("H132").Value = 301283.38
$ws.Range("I132").Value = 515961.47
$ws.Range("K132").Value = 1547884.41
$ws.Range("M132").Value = -1545354.41

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3050.2144
$ws.Range("I136").Value = 985.8
$ws.Range("J136").Value = 4197.1113
$ws.Range("K136").Value = 2957.4
$ws.Range("L136").Value = 12591.3339
$ws.Range("M136").Value = -407.3999999999996
$ws.Range("N136").Value = -17691.3339
